# Insert 2 new weekly observation rows (Pintón / Primera Pintón) right
# before the current row 621, pushing the historical rows down by two
# positions (621-661 -> 623-663).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("621:622").Insert()

# The newly inserted rows 621/622 are blank. Their categorical columns
# (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are identical to the rows that used to sit
# there (now shifted down to 623/624), so copy that shape over first and
# then overwrite just the observation-specific columns (D,M,N,O,P,S).
$ws.Rows.Item(623).Copy()
$ws.Rows.Item(621).PasteSpecial()
$ws.Rows.Item(624).Copy()
$ws.Rows.Item(622).PasteSpecial()
$excel.CutCopyMode = 0

# Row 621: Pintón, week of 2022-06-02 (serial 44714)
$ws.Cells.Item(621, 4).Value = 44714
$ws.Cells.Item(621, 13).Value = 1000
$ws.Cells.Item(621, 14).Value = 10000
$ws.Cells.Item(621, 15).Value = 10000
$ws.Cells.Item(621, 16).Value = 10000
$ws.Cells.Item(621, 19).Value = 500

# Row 622: Primera Pintón, week of 2022-06-02 (serial 44714)
$ws.Cells.Item(622, 4).Value = 44714
$ws.Cells.Item(622, 13).Value = 840
$ws.Cells.Item(622, 14).Value = 11000
$ws.Cells.Item(622, 15).Value = 11000
$ws.Cells.Item(622, 16).Value = 11000
$ws.Cells.Item(622, 19).Value = 550
